$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New "Dynamic Properties" section appended after the existing
# "Upload Download" block (rows 68-73 were blank placeholder rows).

# Row 68: section header
$ws.Range("A68").Value = "dynamicProperties"

# Row 69: nav item
$ws.Range("A69").Value = "dynamicPropertiesNav"
$ws.Range("B69").Value = '//*[@id="item-8"]'
$ws.Range("C69").Value = "By.xpath"

# Row 70: scroll anchor
$ws.Range("A70").Value = "dynamicPropertiesScroll"
$ws.Range("B70").Value = '//*[@id="app"]/div/div/div/div[2]/div[2]/h1'
$ws.Range("C70").Value = "By.xpath"

# Row 71: enable button
$ws.Range("A71").Value = "enableBtn"
$ws.Range("B71").Value = '//*[@id="enableAfter"]'
$ws.Range("C71").Value = "By.xpath"

# Row 72: color button
$ws.Range("A72").Value = "ColorBtn"
$ws.Range("B72").Value = '//*[@id="colorChange"]'
$ws.Range("C72").Value = "By.xpath"

# Row 73: visible button
$ws.Range("A73").Value = "Visiblbtn"
$ws.Range("B73").Value = '//*[@id="visibleAfter"]'
$ws.Range("C73").Value = "By.xpath"

# Match the author's scrolled/selected view state (topLeftCell A61, active cell B73)
$excel.ActiveWindow.ScrollRow = 61
$ws.Range("B73").Select() | Out-Null
